# Edit for CVX_Kenya/Interventions/Shading_nets.xlsx
# "main" sheet (sheet6.xml) - update sensitivity module with new notes/sources
# and refresh a couple of assumption values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")
$ws.Activate()

# --- Row 2: source reference link (replaces old placeholder "BrazAfric") ---
$ws.Range("G2").Value = "https://graduatefarmer.co.ke/marketplace/product/farm-shade-netting/"

# --- Row 5: Cost exponent value updated, and a footnote marker added ---
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "0.4-0.8**"

# --- Row 10: footnote marker appended to productivity range ---
$ws.Range("D10").Value = "0.32-0.99***"

# --- New footnote / references block below the table ---
$ws.Range("A20").Value = "*Based on the climate condition of Kenya for coffee nurseries a shading net with around 70-80% shading level is recommended at the beginning of growing the plants but it can be reduced by the passage of time. (https://www.infonet-biovision.org/PlantHealth/Crops/Coffee)"
$ws.Range("A21").Value = "**Since we have the assumption of purchasing the nets. There is no physical economy of scale"
$ws.Range("A23").Value = "*** We should find an expert to find a more narrow range for productivity based on the change in some factors such as the gowth in photosynthesis and leaves area"
$ws.Range("A24").Value = "ref1: Physiological Growth Response in Seedlings of Arabica Coffee Genotypes Under Contrasting Nursery Microenvironments"
$ws.Range("A25").Value = "ref2: https://www.researchgate.net/publication/293009127_Physiological_Growth_Response_in_Seedlings_of_Arabica_Coffee_Genotypes_Under_Contrasting_Nursery_Microenvironments"

# --- Update the view so the new notes are visible / selected ---
$ws.Range("A21").Select()
